$d = $word.ActiveDocument

# --- Page size & margins ---
$ps = $d.Sections(1).PageSetup
$ps.PageWidth = 419.55
$ps.PageHeight = 595.3
$ps.TopMargin = 49.6
$ps.BottomMargin = 49.6
$ps.LeftMargin = 49.6
$ps.RightMargin = 49.6

# --- Paragraph alignment: left -> both (justify) ---
$d.Paragraphs.Alignment = 3

# --- Font: Times New Roman 10pt -> Book Antiqua 12pt (ascii only, exclude paragraph mark) ---
foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    [void]$r.MoveEnd(1, -1)
    $r.Font.NameAscii = "Book Antiqua"
    $r.Font.Size = 12
}
